$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for refreshed crypto data.
# Leading apostrophe forces text entry so values like '1.00' or
# '69.435.94' are stored as text (matching the sheet's existing inlineStr cells)
# instead of being auto-coerced to numbers by Excel.

$ws.Cells.Item(2, 4).Value = "'69.435.94"
$ws.Cells.Item(2, 5).Value = "'  +1.66%  "
$ws.Cells.Item(3, 4).Value = "'3.389.87"
$ws.Cells.Item(3, 5).Value = "'  +1.13%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "'  +0.11%  "
$ws.Cells.Item(5, 4).Value = "'581.69"
$ws.Cells.Item(5, 5).Value = "'  -0.34%  "
$ws.Cells.Item(6, 4).Value = "'179.04"
$ws.Cells.Item(6, 5).Value = "'  +0.86%  "
$ws.Cells.Item(7, 5).Value = "'  +0.07%  "
$ws.Cells.Item(8, 5).Value = "'  +0.56%  "
$ws.Cells.Item(9, 5).Value = "'  +7.92%  "
$ws.Cells.Item(10, 5).Value = "'  +0.82%  "
$ws.Cells.Item(11, 4).Value = "'48.40"
$ws.Cells.Item(11, 5).Value = "'  +0.24%  "
$ws.Cells.Item(12, 5).Value = "'  +3.60%  "
$ws.Cells.Item(13, 4).Value = "'685.72"
$ws.Cells.Item(13, 5).Value = "'  +0.00%  "
$ws.Cells.Item(14, 4).Value = "'8.60"
$ws.Cells.Item(14, 5).Value = "'  +2.07%  "
$ws.Cells.Item(15, 4).Value = "'3.934.63"
$ws.Cells.Item(15, 5).Value = "'  +1.02%  "
$ws.Cells.Item(16, 4).Value = "'69.546.82"
$ws.Cells.Item(16, 5).Value = "'  +1.78%  "
$ws.Cells.Item(17, 4).Value = "'0.120"
$ws.Cells.Item(17, 5).Value = "'  +0.80%  "
$ws.Cells.Item(18, 4).Value = "'3.385.95"
$ws.Cells.Item(18, 5).Value = "'  +0.96%  "
$ws.Cells.Item(19, 4).Value = "'17.70"
$ws.Cells.Item(19, 5).Value = "'  +1.38%  "
$ws.Cells.Item(20, 4).Value = "'11.29"
$ws.Cells.Item(20, 5).Value = "'  +0.80%  "
$ws.Cells.Item(21, 4).Value = "'0.910"
$ws.Cells.Item(21, 5).Value = "'  +1.64%  "
$ws.Cells.Item(22, 4).Value = "'17.16"
$ws.Cells.Item(23, 4).Value = "'5.36"
$ws.Cells.Item(23, 5).Value = "'  -1.44%  "
$ws.Cells.Item(24, 4).Value = "'101.20"
$ws.Cells.Item(24, 5).Value = "'  +1.29%  "
$ws.Cells.Item(25, 4).Value = "'3.89"
$ws.Cells.Item(25, 5).Value = "'  -0.73%  "
$ws.Cells.Item(26, 5).Value = "'  +0.08%  "
$ws.Cells.Item(27, 4).Value = "'9.75"
$ws.Cells.Item(27, 5).Value = "'  +2.05%  "
$ws.Cells.Item(28, 4).Value = "'33.46"
$ws.Cells.Item(28, 5).Value = "'  +1.66%  "
$ws.Cells.Item(29, 5).Value = "'  +2.70%  "
$ws.Cells.Item(30, 5).Value = "'  +0.23%  "
$ws.Cells.Item(31, 4).Value = "'3.85"
$ws.Cells.Item(31, 5).Value = "'  +17.11%  "
$ws.Cells.Item(32, 4).Value = "'11.04"
$ws.Cells.Item(33, 4).Value = "'550.26"
$ws.Cells.Item(33, 5).Value = "'  -2.14%  "
$ws.Cells.Item(34, 5).Value = "'  +0.13%  "
$ws.Cells.Item(35, 4).Value = "'57.89"
$ws.Cells.Item(35, 5).Value = "'  -0.01%  "
$ws.Cells.Item(36, 5).Value = "'  +0.16%  "
$ws.Cells.Item(37, 4).Value = "'3.606.10"
$ws.Cells.Item(37, 5).Value = "'  -2.64%  "
$ws.Cells.Item(38, 5).Value = "'  +3.21%  "
$ws.Cells.Item(39, 4).Value = "'35.49"
$ws.Cells.Item(39, 5).Value = "'  +2.48%  "
$ws.Cells.Item(40, 4).Value = "'0.0₃0743"
$ws.Cells.Item(40, 5).Value = "'  +10.43%  "
$ws.Cells.Item(41, 5).Value = "'  +5.16%  "
$ws.Cells.Item(42, 4).Value = "'2.73"
$ws.Cells.Item(42, 5).Value = "'  +4.63%  "
$ws.Cells.Item(43, 5).Value = "'  +3.65%  "
$ws.Cells.Item(44, 4).Value = "'0.0426"
$ws.Cells.Item(44, 5).Value = "'  +3.46%  "
$ws.Cells.Item(45, 5).Value = "'  +0.04%  "
$ws.Cells.Item(46, 5).Value = "'  +0.85%  "
$ws.Cells.Item(47, 5).Value = "'  +0.27%  "
$ws.Cells.Item(48, 5).Value = "'  +3.73%  "
$ws.Cells.Item(49, 4).Value = "'1.00"
$ws.Cells.Item(49, 5).Value = "'  -0.17%  "
$ws.Cells.Item(50, 4).Value = "'129.78"
$ws.Cells.Item(50, 5).Value = "'  -0.97%  "
$ws.Cells.Item(51, 4).Value = "'2.58"
$ws.Cells.Item(51, 5).Value = "'  +1.10%  "
